$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "day" sheet: D431..D434 were stored as text ("532454" etc.) - convert
#    them to genuine numeric cells (same displayed value).
# ---------------------------------------------------------------------------
$day = $wb.Worksheets.Item("day")
$day.Range("D431").Value = 532454
$day.Range("D432").Value = 532187
$day.Range("D433").Value = 531213
$day.Range("D434").Value = 500469

# ---------------------------------------------------------------------------
# 2) "week" sheet: append 12 new rows (158-169) of stock data.
# ---------------------------------------------------------------------------
$week = $wb.Worksheets.Item("week")

function Add-WeekRow {
    param(
        [int]$Row,
        [int]$Sr,
        [string]$NseCode,
        [string]$Name,
        [string]$BseCode,
        [double]$PerChg,
        [double]$Close,
        [double]$Volume,
        [string]$Timeframe,
        [string]$DateTime
    )

    $week.Cells.Item($Row, 1).Value = $Sr
    $week.Cells.Item($Row, 2).Value = $NseCode
    $week.Cells.Item($Row, 3).Value = $Name

    # bsecode must stay a text cell even though it looks numeric - format the
    # cell as Text before writing, then drop the formatting again so no
    # leftover style index is attached to the cell.
    $dcell = $week.Cells.Item($Row, 4)
    $dcell.NumberFormat = "@"
    $dcell.Value = $BseCode
    $dcell.ClearFormats()

    $week.Cells.Item($Row, 5).Value = $PerChg
    $week.Cells.Item($Row, 6).Value = $Close
    $week.Cells.Item($Row, 7).Value = $Volume
    $week.Cells.Item($Row, 8).Value = $Timeframe
    $week.Cells.Item($Row, 9).Value = $DateTime
}

Add-WeekRow 158 1  "SHREECEM"   "Shree Cements Limited"              "500387" -1.22     24706.05 28827    "week" "23/08/2024 11:33:20"
Add-WeekRow 159 2  "LTIM"       "LTI Mindtree Ltd"                   "540005" -1.1      5641.6   252213   "week" "23/08/2024 11:33:20"
Add-WeekRow 160 3  "ASIANPAINT" "Asian Paints Limited"                "500820" -1        3154.65  779284   "week" "23/08/2024 11:33:20"
Add-WeekRow 161 4  "GODREJPROP" "Godrej Properties Limited"           "533150" -1.69     2883.65  496462   "week" "23/08/2024 11:33:20"
Add-WeekRow 162 5  "DALBHARAT"  "Dalmia Bharat Limited"                "533309" 0.26      1804.05  311355   "week" "23/08/2024 11:33:20"
Add-WeekRow 163 6  "OBEROIRLTY" "Oberoi Realty Limited"                "533273" -1.79     1703.1   570019   "week" "23/08/2024 11:33:20"
Add-WeekRow 164 7  "SBIN"       "State Bank Of India"                  "500112" -0.6      815.35   5537747  "week" "23/08/2024 11:33:20"
Add-WeekRow 165 8  "LICHSGFIN"  "Lic Housing Finance Limited"          "500253" -1.43     672.6    1556950  "week" "23/08/2024 11:33:20"
Add-WeekRow 166 9  "INDHOTEL"   "The Indian Hotels Company Limited"    "500850" -0.42     641.9    3388830  "week" "23/08/2024 11:33:20"
Add-WeekRow 167 10 "BSOFT"      "Birlasoft Ltd"                        "532400" -0.07     601.25   3030044  "week" "23/08/2024 11:33:20"
Add-WeekRow 168 11 "WIPRO"      "Wipro Limited"                        "507685" -1.27     512.4    5638853  "week" "23/08/2024 11:33:20"
Add-WeekRow 169 12 "MANAPPURAM" "Manappuram Finance Limited"           "531213" -0.79     215.08   6489947  "week" "23/08/2024 11:33:20"
